$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells keep trailing zeros / punctuation as literal text
# rather than being reinterpreted as numbers when assigned via .Value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.016.93"
$ws.Range("E2").Value = "  -2.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.81"
$ws.Range("E3").Value = "  -1.92%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.54"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5098"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2629"
$ws.Range("E8").Value = "  -0.88%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06390"
$ws.Range("E9").Value = "  +2.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.75"
$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07416"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.666.27"
$ws.Range("E12").Value = "  -1.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.500"
$ws.Range("E13").Value = "  +0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5809"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008512"
$ws.Range("E15").Value = "  +1.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.19"
$ws.Range("E16").Value = "  -2.14%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.086.54"
$ws.Range("E17").Value = "  -2.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.900"
$ws.Range("E18").Value = "  -3.07%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.68"
$ws.Range("E20").Value = "  -1.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.68"
$ws.Range("E21").Value = "  +0.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.200"
$ws.Range("E22").Value = "  -1.38%  "

$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.62"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.594"
$ws.Range("E25").Value = "  -0.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1188"
$ws.Range("E26").Value = "  +3.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.61"
$ws.Range("E27").Value = "  -0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06720"
$ws.Range("E28").Value = "  +17.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.304"
$ws.Range("E29").Value = "  +0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.318"
$ws.Range("E30").Value = "  -1.34%  "

$ws.Range("E31").Value = "  +0.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.502"
$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("E33").Value = "  -2.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.017"
$ws.Range("E34").Value = "  -0.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6060"
$ws.Range("E35").Value = "  +0.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.365"
$ws.Range("E36").Value = "  -0.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.680"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.209"
$ws.Range("E38").Value = "  +5.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01608"
$ws.Range("E39").Value = "  +0.64%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.074.81"
$ws.Range("E40").Value = "  -1.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8597"
$ws.Range("E41").Value = "  -0.74%  "

$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("E43").Value = "  +1.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.24"
$ws.Range("E46").Value = "  -1.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.022"
$ws.Range("E48").Value = "  -1.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05211"
$ws.Range("E49").Value = "  -0.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4288"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.944"
$ws.Range("E51").Value = "  +3.06%  "

# Rows 44-45: BabyDogeCoin and RocketPoolETH swap positions with refreshed values
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.00000000116"
$ws.Range("E44").Value = "  +8.49%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.813.60"
$ws.Range("E45").Value = "  -2.25%  "
